# Empfaenger-Upload-Links-7.xlsx — "Add files via upload"
#
# Rename the header in B1 from "Link" to "UploadLink", and move the
# worksheet's active cell/selection from C1 to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 header: "Link" -> "UploadLink"
$ws.Range("B1").Value = "UploadLink"

# Update the current selection to B2 (was C1)
$ws.Range("B2").Select()
